# Re-run of the NATMI ligand-receptor analysis for Fgf2-Fgfr2 with updated TPM
# values: the "Resolving-Mac" cluster is no longer among the recipient (target)
# clusters, shrinking the data block from 4 clusters x 4 clusters (16 rows) to
# 4 clusters x 3 target clusters (12 rows), and every metric column recomputed
# for the remaining sending/target cluster pairs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    @("ECs", "Fgf2", "Fgfr2", "ECs", 1, 0.3333333333333333, 0.313179, 0.939537, 0.02707464596575709, 0.0270746459657571, 3, 1, 1.414593, 4.243779, 0.3478871232761722, 0.3478871232761722, 0.443020821147, 3.987187390323, 0.009418920698748055, 0.009418920698748058),
    @("ECs", "Fgf2", "Fgfr2", "FAPs", 1, 0.3333333333333333, 0.313179, 0.939537, 0.02707464596575709, 0.0270746459657571, 3, 1, 2.337487333333333, 7.012461999999999, 0.5748520910875596, 0.5748520910875596, 0.732051945566, 6.588467510093999, 0.01556391684887083, 0.01556391684887083),
    @("ECs", "Fgf2", "Fgfr2", "MuSCs", 1, 0.3333333333333333, 0.313179, 0.939537, 0.02707464596575709, 0.0270746459657571, 3, 1, 0.314161, 0.942483, 0.07726078563626818, 0.07726078563626819, 0.09838862781899999, 0.8854976503709999, 0.002091808418138212, 0.002091808418138213),
    @("FAPs", "Fgf2", "Fgfr2", "ECs", 3, 1, 8.077181333333334, 24.231544, 0.6982806158817221, 0.6982806158817222, 3, 1, 1.414593, 4.243779, 0.3478871232761722, 0.3478871232761722, 11.425924173864, 102.833317564776, 0.2429228346986061, 0.2429228346986061),
    @("FAPs", "Fgf2", "Fgfr2", "FAPs", 3, 1, 8.077181333333334, 24.231544, 0.6982806158817221, 0.6982806158817222, 3, 1, 2.337487333333333, 7.012461999999999, 0.5748520910875596, 0.5748520910875596, 18.88030905570311, 169.922781501328, 0.4014080722055169, 0.401408072205517),
    @("FAPs", "Fgf2", "Fgfr2", "MuSCs", 3, 1, 8.077181333333334, 24.231544, 0.6982806158817221, 0.6982806158817222, 3, 1, 0.314161, 0.942483, 0.07726078563626818, 0.07726078563626819, 2.537535364861333, 22.837818283752, 0.05394970897759905, 0.05394970897759906),
    @("MuSCs", "Fgf2", "Fgfr2", "ECs", 3, 1, 2.897745666666667, 8.693237, 0.2505130868410934, 0.2505130868410934, 3, 1, 1.414593, 4.243779, 0.3478871232761722, 0.3478871232761722, 4.099130735847, 36.892176622623, 0.08715027712418187, 0.08715027712418189),
    @("MuSCs", "Fgf2", "Fgfr2", "FAPs", 3, 1, 2.897745666666667, 8.693237, 0.2505130868410934, 0.2505130868410934, 3, 1, 2.337487333333333, 7.012461999999999, 0.5748520910875596, 0.5748520910875596, 6.773443791054889, 60.960994119494, 0.1440079718154019, 0.1440079718154019),
    @("MuSCs", "Fgf2", "Fgfr2", "MuSCs", 3, 1, 2.897745666666667, 8.693237, 0.2505130868410934, 0.2505130868410934, 3, 1, 0.314161, 0.942483, 0.07726078563626818, 0.07726078563626819, 0.9103586763856666, 8.193228087470999, 0.01935483790150955, 0.01935483790150955),
    @("Resolving-Mac", "Fgf2", "Fgfr2", "ECs", 2, 0.6666666666666666, 0.2791366666666666, 0.83741, 0.02413165131142748, 0.02413165131142749, 3, 1, 1.414593, 4.243779, 0.3478871232761722, 0.3478871232761722, 0.39486477471, 3.55378297239, 0.008395090754636176, 0.008395090754636179),
    @("Resolving-Mac", "Fgf2", "Fgfr2", "FAPs", 2, 0.6666666666666666, 0.2791366666666666, 0.83741, 0.02413165131142748, 0.02413165131142749, 3, 1, 2.337487333333333, 7.012461999999999, 0.5748520910875596, 0.5748520910875596, 0.6524784226022221, 5.87230580342, 0.01387213021776994, 0.01387213021776994),
    @("Resolving-Mac", "Fgf2", "Fgfr2", "MuSCs", 2, 0.6666666666666666, 0.2791366666666666, 0.83741, 0.02413165131142748, 0.02413165131142749, 3, 1, 0.314161, 0.942483, 0.07726078563626818, 0.07726078563626819, 0.08769385433666665, 0.78924468903, 0.001864430339021369, 0.00186443033902137)
)

for ($i = 0; $i -lt $rowData.Count; $i++) {
    $r = $i + 2
    $rowVals = $rowData[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}

# Remove rows 14-17, which no longer have data in the updated dataset
$ws.Range("A14:T17").Delete() | Out-Null
